$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "56.67"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "2.89"
$ws.Range("E45").Value = "  +4.23%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "32.02"
$ws.Range("E46").Value = "  -1.39%  "

$ws.Range("D2").Value = "69.697.30"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "3.493.92"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "599.86"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").Value = "170.77"
$ws.Range("E6").Value = "  +2.44%  "
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").Value = "3.493.89"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.191"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").Value = "7.27"
$ws.Range("E11").Value = "  +6.18%  "
$ws.Range("D12").Value = "0.577"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "45.79"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "0.0000272"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "4.063.21"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "8.21"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "600.99"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "3.508.54"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").Value = "69.845.07"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").Value = "17.04"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "0.863"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "9.17"
$ws.Range("E23").Value = "  -16.52%  "
$ws.Range("D24").Value = "15.47"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").Value = "95.12"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "3.69"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "2.55"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "33.68"
$ws.Range("E29").Value = "  +3.74%  "
$ws.Range("D30").Value = "8.90"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "703.18"
$ws.Range("E31").Value = "  +21.68%  "
$ws.Range("D32").Value = "2.99"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").Value = "8.04"
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("D34").Value = "6.86"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").Value = "0.0990"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").Value = "3.53"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").Value = "10.62"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").Value = "0.0469"
$ws.Range("E39").Value = "  +8.21%  "
$ws.Range("D42").Value = "0.141"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").Value = "3.321.57"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").Value = "0.312"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("D47").Value = "0.0₃0682"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "2.52"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "0.129"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "132.15"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  -0.02%  "
